$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the source values (A2, A3); dependent formulas (C2, D2, C3, D3, D6, E6)
# recalculate automatically.
$ws.Range("A2").Value = 65
$ws.Range("A3").Value = 140

# Update the selected range shown in the sheet view.
$ws.Range("I21:I22").Select() | Out-Null
